$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per the diff
$ws.Range("D2").Value = '51.937.03'
$ws.Range("E2").Value = '  +0.00%  '

$ws.Range("D3").Value = '2.774.68'
$ws.Range("E3").Value = '  -1.67%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").Value = '355.82'
$ws.Range("E5").Value = '  +0.09%  '

$ws.Range("D6").Value = '108.83'
$ws.Range("E6").Value = '  -4.23%  '

$ws.Range("E7").Value = '  +1.63%  '

$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").Value = '0.589'
$ws.Range("E9").Value = '  -1.86%  '

$ws.Range("D10").Value = '40.15'
$ws.Range("E10").Value = '  -4.43%  '

$ws.Range("D11").Value = '0.0851'
$ws.Range("E11").Value = '  +0.01%  '

$ws.Range("E12").Value = '  +0.78%  '

$ws.Range("D13").Value = '19.36'
$ws.Range("E13").Value = '  -3.45%  '

$ws.Range("D14").Value = '7.59'
$ws.Range("E14").Value = '  -1.68%  '

$ws.Range("D15").Value = '3.214.65'
$ws.Range("E15").Value = '  -0.98%  '

$ws.Range("D16").Value = '2.788.55'
$ws.Range("E16").Value = '  -1.58%  '

$ws.Range("D17").Value = '0.932'
$ws.Range("E17").Value = '  +3.86%  '

$ws.Range("D18").Value = '51.770.87'
$ws.Range("E18").Value = '  -0.15%  '

$ws.Range("D19").Value = '7.40'
$ws.Range("E19").Value = '  +0.06%  '

$ws.Range("E20").Value = '  -1.28%  '

$ws.Range("D21").Value = '13.01'
$ws.Range("E21").Value = '  -4.16%  '

$ws.Range("D22").Value = '0.0₃0975'
$ws.Range("E22").Value = '  -2.31%  '

$ws.Range("D23").Value = '274.16'
$ws.Range("E23").Value = '  +1.55%  '

$ws.Range("D24").Value = '69.80'
$ws.Range("E24").Value = '  +0.10%  '

$ws.Range("D25").Value = '2.72'
$ws.Range("E25").Value = '  -2.52%  '

$ws.Range("D26").Value = '26.50'
$ws.Range("E26").Value = '  -1.19%  '

$ws.Range("E27").Value = '  -0.01%  '

$ws.Range("D28").Value = '10.12'
$ws.Range("E28").Value = '  -1.91%  '

$ws.Range("D30").Value = '0.143'
$ws.Range("E30").Value = '  +2.05%  '

$ws.Range("D31").Value = '0.0465'
$ws.Range("E31").Value = '  +2.88%  '

$ws.Range("D32").Value = '51.52'
$ws.Range("E32").Value = '  +1.39%  '

$ws.Range("E33").Value = '  -0.54%  '

$ws.Range("D34").Value = '5.69'
$ws.Range("E34").Value = '  -2.50%  '

$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D35").Value = '5.32'
$ws.Range("E35").Value = '  +8.80%  '

$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = '0.0842'
$ws.Range("E36").Value = '  +1.32%  '

$ws.Range("E37").Value = '  +0.07%  '

$ws.Range("D38").Value = '3.20'
$ws.Range("E38").Value = '  -0.45%  '

$ws.Range("D39").Value = '18.03'
$ws.Range("E39").Value = '  -2.03%  '

$ws.Range("E40").Value = '  -4.87%  '

$ws.Range("E41").Value = '  -1.89%  '

$ws.Range("D42").Value = '0.115'
$ws.Range("E42").Value = '  -0.40%  '

$ws.Range("E43").Value = '  -2.89%  '

$ws.Range("D44").Value = '121.30'
$ws.Range("E44").Value = '  -5.51%  '

$ws.Range("D45").Value = '21.86'
$ws.Range("E45").Value = '  -7.81%  '

$ws.Range("D46").Value = '2.052.17'
$ws.Range("E46").Value = '  -1.14%  '

$ws.Range("E47").Value = '  -3.18%  '

$ws.Range("E48").Value = '  -1.97%  '

$ws.Range("D49").Value = '5.68'
$ws.Range("E49").Value = '  -0.23%  '

$ws.Range("E50").Value = '  -3.94%  '

$ws.Range("D51").Value = '8.89'
$ws.Range("E51").Value = '  -0.29%  '
